$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 91.94136433333334
$ws.Range("H2").Value = 275.824093
$ws.Range("I2").Value = 0.02307547609860541
$ws.Range("J2").Value = 0.02307547609860541
$ws.Range("M2").Value = 211.2725676666666
$ws.Range("N2").Value = 633.8177029999999
$ws.Range("O2").Value = 0.6324644927232657
$ws.Range("P2").Value = 0.6324644927232657
$ws.Range("Q2").Value = 19424.68811747982
$ws.Range("R2").Value = 174822.1930573184
$ws.Range("S2").Value = 0.01459441928505231
$ws.Range("T2").Value = 0.01459441928505231

$ws.Range("G3").Value = 91.94136433333334
$ws.Range("H3").Value = 275.824093
$ws.Range("I3").Value = 0.02307547609860541
$ws.Range("J3").Value = 0.02307547609860541
$ws.Range("M3").Value = 59.36675400000001
$ws.Range("N3").Value = 178.100262
$ws.Range("O3").Value = 0.1777200152765546
$ws.Range("P3").Value = 0.1777200152765546
$ws.Range("Q3").Value = 5458.260358801375
$ws.Range("R3").Value = 49124.34322921237
$ws.Range("S3").Value = 0.004100973964757925
$ws.Range("T3").Value = 0.004100973964757924

$ws.Range("G4").Value = 91.94136433333334
$ws.Range("H4").Value = 275.824093
$ws.Range("I4").Value = 0.02307547609860541
$ws.Range("J4").Value = 0.02307547609860541
$ws.Range("M4").Value = 0.4593846666666666
$ws.Range("N4").Value = 1.378154
$ws.Range("O4").Value = 0.001375211620595172
$ws.Range("P4").Value = 0.001375211620595172
$ws.Range("Q4").Value = 42.23645300714688
$ws.Range("R4").Value = 380.128077064322
$ws.Range("S4").Value = 0.00003173366288156831
$ws.Range("T4").Value = 0.00003173366288156831

$ws.Range("G5").Value = 91.94136433333334
$ws.Range("H5").Value = 275.824093
$ws.Range("I5").Value = 0.02307547609860541
$ws.Range("J5").Value = 0.02307547609860541
$ws.Range("M5").Value = 62.94782133333333
$ws.Range("N5").Value = 188.843464
$ws.Range("O5").Value = 0.1884402803795846
$ws.Range("P5").Value = 0.1884402803795846
$ws.Range("Q5").Value = 5787.508575197572
$ws.Range("R5").Value = 52087.57717677815
$ws.Range("S5").Value = 0.004348349185913608
$ws.Range("T5").Value = 0.004348349185913607

$ws.Range("I6").Value = 0.9681738695089209
$ws.Range("J6").Value = 0.9681738695089209
$ws.Range("M6").Value = 211.2725676666666
$ws.Range("N6").Value = 633.8177029999999
$ws.Range("O6").Value = 0.6324644927232657
$ws.Range("P6").Value = 0.6324644927232657
$ws.Range("Q6").Value = 814998.3722260438
$ws.Range("R6").Value = 7334985.350034394
$ws.Range("S6").Value = 0.6123355952468809
$ws.Range("T6").Value = 0.6123355952468809

$ws.Range("I7").Value = 0.9681738695089209
$ws.Range("J7").Value = 0.9681738695089209
$ws.Range("M7").Value = 59.36675400000001
$ws.Range("N7").Value = 178.100262
$ws.Range("O7").Value = 0.1777200152765546
$ws.Range("P7").Value = 0.1777200152765546
$ws.Range("Q7").Value = 229011.3118267256
$ws.Range("R7").Value = 2061101.80644053
$ws.Range("S7").Value = 0.1720638748794864
$ws.Range("T7").Value = 0.1720638748794864

$ws.Range("I8").Value = 0.9681738695089209
$ws.Range("J8").Value = 0.9681738695089209
$ws.Range("M8").Value = 0.4593846666666666
$ws.Range("N8").Value = 1.378154
$ws.Range("O8").Value = 0.001375211620595172
$ws.Range("P8").Value = 0.001375211620595172
$ws.Range("Q8").Value = 1772.107754896223
$ws.Range("R8").Value = 15948.96979406601
$ws.Range("S8").Value = 0.001331443956105262
$ws.Range("T8").Value = 0.001331443956105262

$ws.Range("I9").Value = 0.9681738695089209
$ws.Range("J9").Value = 0.9681738695089209
$ws.Range("M9").Value = 62.94782133333333
$ws.Range("N9").Value = 188.843464
$ws.Range("O9").Value = 0.1884402803795846
$ws.Range("P9").Value = 0.1884402803795846
$ws.Range("Q9").Value = 242825.5238644344
$ws.Range("R9").Value = 2185429.71477991
$ws.Range("S9").Value = 0.1824429554264484
$ws.Range("T9").Value = 0.1824429554264484

$ws.Range("G10").Value = 1.840730666666667
$ws.Range("H10").Value = 5.522192
$ws.Range("I10").Value = 0.000461987232956876
$ws.Range("J10").Value = 0.000461987232956876
$ws.Range("M10").Value = 211.2725676666666
$ws.Range("N10").Value = 633.8177029999999
$ws.Range("O10").Value = 0.6324644927232657
$ws.Range("P10").Value = 0.6324644927232657
$ws.Range("Q10").Value = 388.8958943294418
$ws.Range("R10").Value = 3500.063048964976
$ws.Range("S10").Value = 0.0002921905209366958
$ws.Range("T10").Value = 0.0002921905209366958

$ws.Range("G11").Value = 1.840730666666667
$ws.Range("H11").Value = 5.522192
$ws.Range("I11").Value = 0.000461987232956876
$ws.Range("J11").Value = 0.000461987232956876
$ws.Range("M11").Value = 59.36675400000001
$ws.Range("N11").Value = 178.100262
$ws.Range("O11").Value = 0.1777200152765546
$ws.Range("P11").Value = 0.1777200152765546
$ws.Range("Q11").Value = 109.278204668256
$ws.Range("R11").Value = 983.5038420143042
$ws.Range("S11").Value = 0.00008210437809866918
$ws.Range("T11").Value = 0.00008210437809866918

$ws.Range("G12").Value = 1.840730666666667
$ws.Range("H12").Value = 5.522192
$ws.Range("I12").Value = 0.000461987232956876
$ws.Range("J12").Value = 0.000461987232956876
$ws.Range("M12").Value = 0.4593846666666666
$ws.Range("N12").Value = 1.378154
$ws.Range("O12").Value = 0.001375211620595172
$ws.Range("P12").Value = 0.001375211620595172
$ws.Range("Q12").Value = 0.8456034437297777
$ws.Range("R12").Value = 7.610430993568
$ws.Range("S12").Value = 0.0000006353302113289048
$ws.Range("T12").Value = 0.0000006353302113289048

$ws.Range("G13").Value = 1.840730666666667
$ws.Range("H13").Value = 5.522192
$ws.Range("I13").Value = 0.000461987232956876
$ws.Range("J13").Value = 0.000461987232956876
$ws.Range("M13").Value = 62.94782133333333
$ws.Range("N13").Value = 188.843464
$ws.Range("O13").Value = 0.1884402803795846
$ws.Range("P13").Value = 0.1884402803795846
$ws.Range("Q13").Value = 115.8699851281209
$ws.Range("R13").Value = 1042.829866153088
$ws.Range("S13").Value = 0.00008705700371018219
$ws.Range("T13").Value = 0.00008705700371018219

$ws.Range("G14").Value = 31.11921133333334
$ws.Range("H14").Value = 93.357634
$ws.Range("I14").Value = 0.007810310653280575
$ws.Range("J14").Value = 0.007810310653280575
$ws.Range("M14").Value = 211.2725676666666
$ws.Range("N14").Value = 633.8177029999999
$ws.Range("O14").Value = 0.6324644927232657
$ws.Range("P14").Value = 0.6324644927232657
$ws.Range("Q14").Value = 6574.635682154967
$ws.Range("R14").Value = 59171.7211393947
$ws.Range("S14").Value = 0.004939744165338216
$ws.Range("T14").Value = 0.004939744165338216

$ws.Range("G15").Value = 31.11921133333334
$ws.Range("H15").Value = 93.357634
$ws.Range("I15").Value = 0.007810310653280575
$ws.Range("J15").Value = 0.007810310653280575
$ws.Range("M15").Value = 59.36675400000001
$ws.Range("N15").Value = 178.100262
$ws.Range("O15").Value = 0.1777200152765546
$ws.Range("P15").Value = 0.1777200152765546
$ws.Range("Q15").Value = 1847.446563900012
$ws.Range("R15").Value = 16627.01907510011
$ws.Range("S15").Value = 0.001388048528615661
$ws.Range("T15").Value = 0.001388048528615661

$ws.Range("G16").Value = 31.11921133333334
$ws.Range("H16").Value = 93.357634
$ws.Range("I16").Value = 0.007810310653280575
$ws.Range("J16").Value = 0.007810310653280575
$ws.Range("M16").Value = 0.4593846666666666
$ws.Range("N16").Value = 1.378154
$ws.Range("O16").Value = 0.001375211620595172
$ws.Range("P16").Value = 0.001375211620595172
$ws.Range("Q16").Value = 14.29568852529289
$ws.Range("R16").Value = 128.661196727636
$ws.Range("S16").Value = 0.00001074082997084972
$ws.Range("T16").Value = 0.00001074082997084972

$ws.Range("G17").Value = 31.11921133333334
$ws.Range("H17").Value = 93.357634
$ws.Range("I17").Value = 0.007810310653280575
$ws.Range("J17").Value = 0.007810310653280575
$ws.Range("M17").Value = 62.94782133333333
$ws.Range("N17").Value = 188.843464
$ws.Range("O17").Value = 0.1884402803795846
$ws.Range("P17").Value = 0.1884402803795846
$ws.Range("Q17").Value = 1958.886555044909
$ws.Range("R17").Value = 17629.97899540417
$ws.Range("S17").Value = 0.001471777129355848
$ws.Range("T17").Value = 0.001471777129355848

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.905952
$ws.Range("H18").Value = 5.717856
$ws.Range("I18").Value = 0.0004783565062362683
$ws.Range("J18").Value = 0.0004783565062362683
$ws.Range("M18").Value = 211.2725676666666
$ws.Range("N18").Value = 633.8177029999999
$ws.Range("O18").Value = 0.6324644927232657
$ws.Range("P18").Value = 0.6324644927232657
$ws.Range("Q18").Value = 402.6753728894187
$ws.Range("R18").Value = 3624.078356004768
$ws.Range("S18").Value = 0.0003025435050575951
$ws.Range("T18").Value = 0.0003025435050575951

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.905952
$ws.Range("H19").Value = 5.717856
$ws.Range("I19").Value = 0.0004783565062362683
$ws.Range("J19").Value = 0.0004783565062362683
$ws.Range("M19").Value = 59.36675400000001
$ws.Range("N19").Value = 178.100262
$ws.Range("O19").Value = 0.1777200152765546
$ws.Range("P19").Value = 0.1777200152765546
$ws.Range("Q19").Value = 113.150183519808
$ws.Range("R19").Value = 1018.351651678272
$ws.Range("S19").Value = 0.00008501352559594889
$ws.Range("T19").Value = 0.00008501352559594888

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1.905952
$ws.Range("H20").Value = 5.717856
$ws.Range("I20").Value = 0.0004783565062362683
$ws.Range("J20").Value = 0.0004783565062362683
$ws.Range("M20").Value = 0.4593846666666666
$ws.Range("N20").Value = 1.378154
$ws.Range("O20").Value = 0.001375211620595172
$ws.Range("P20").Value = 0.001375211620595172
$ws.Range("Q20").Value = 0.8755651242026666
$ws.Range("R20").Value = 7.880086117824
$ws.Range("S20").Value = 0.0000006578414261634232
$ws.Range("T20").Value = 0.0000006578414261634232

$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1.905952
$ws.Range("H21").Value = 5.717856
$ws.Range("I21").Value = 0.0004783565062362683
$ws.Range("J21").Value = 0.0004783565062362683
$ws.Range("M21").Value = 62.94782133333333
$ws.Range("N21").Value = 188.843464
$ws.Range("O21").Value = 0.1884402803795846
$ws.Range("P21").Value = 0.1884402803795846
$ws.Range("Q21").Value = 119.9755259659093
$ws.Range("R21").Value = 1079.779733693184
$ws.Range("S21").Value = 0.00009014163415656093
$ws.Range("T21").Value = 0.00009014163415656091

